$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # C社 実績値input (actuals)
$ws2 = $wb.Worksheets.Item(2)  # C社 予測値input (forecast)

# ---- Sheet1 (actuals): extend with new rows 39-50 (202104-202203) ----
# Copy formatting (fill styles) from row 38 down to the new rows first
$ws1.Range("B38:C38").Copy()
$ws1.Range("B39:C50").PasteSpecial(-4122)

$ws1.Cells.Item(39, 2).Value = 202104
$ws1.Cells.Item(39, 3).Value = 92800
$ws1.Cells.Item(40, 2).Value = 202105
$ws1.Cells.Item(40, 3).Value = 87800
$ws1.Cells.Item(41, 2).Value = 202106
$ws1.Cells.Item(41, 3).Value = 109780
$ws1.Cells.Item(42, 2).Value = 202107
$ws1.Cells.Item(42, 3).Value = 142400
$ws1.Cells.Item(43, 2).Value = 202108
$ws1.Cells.Item(43, 3).Value = 99040
$ws1.Cells.Item(44, 2).Value = 202109
$ws1.Cells.Item(44, 3).Value = 129600
$ws1.Cells.Item(45, 2).Value = 202110
$ws1.Cells.Item(45, 3).Value = 72400
$ws1.Cells.Item(46, 2).Value = 202111
$ws1.Cells.Item(46, 3).Value = 126000
$ws1.Cells.Item(47, 2).Value = 202112
$ws1.Cells.Item(47, 3).Value = 181400
$ws1.Cells.Item(48, 2).Value = 202201
$ws1.Cells.Item(48, 3).Value = 113200
$ws1.Cells.Item(49, 2).Value = 202202
$ws1.Cells.Item(49, 3).Value = 84000
$ws1.Cells.Item(50, 2).Value = 202203
$ws1.Cells.Item(50, 3).Value = 97200

# ---- Sheet2 (forecast): fill previously-unknown cells in rows 35-38 ----
$ws2.Cells.Item(35, 3).Value = 79200
$ws2.Cells.Item(35, 4).Value = 80000
$ws2.Cells.Item(35, 5).Value = 84000
$ws2.Cells.Item(35, 6).Value = 86000
$ws2.Cells.Item(35, 7).Value = 96200
$ws2.Cells.Item(35, 8).Value = 96200
$ws2.Cells.Item(35, 9).Value = 96200
$ws2.Cells.Item(35, 10).Value = 103200
$ws2.Cells.Item(35, 11).Value = 113200
$ws2.Cells.Item(35, 12).Value = 113200
$ws2.Cells.Item(35, 13).Value = 122200
$ws2.Cells.Item(35, 14).Value = 122200

$ws2.Cells.Item(36, 3).Value = 76000
$ws2.Cells.Item(36, 4).Value = 76400
$ws2.Cells.Item(36, 5).Value = 93000
$ws2.Cells.Item(36, 6).Value = 96200
$ws2.Cells.Item(36, 7).Value = 96200
$ws2.Cells.Item(36, 8).Value = 96200
$ws2.Cells.Item(36, 9).Value = 103200
$ws2.Cells.Item(36, 10).Value = 113200
$ws2.Cells.Item(36, 11).Value = 113200
$ws2.Cells.Item(36, 12).Value = 122200
$ws2.Cells.Item(36, 13).Value = 122200
$ws2.Cells.Item(36, 14).Value = 122200

$ws2.Cells.Item(37, 3).Value = 100400
$ws2.Cells.Item(37, 4).Value = 95000
$ws2.Cells.Item(37, 5).Value = 96200
$ws2.Cells.Item(37, 6).Value = 96200
$ws2.Cells.Item(37, 7).Value = 96200
$ws2.Cells.Item(37, 8).Value = 103200
$ws2.Cells.Item(37, 9).Value = 113200
$ws2.Cells.Item(37, 10).Value = 113200
$ws2.Cells.Item(37, 11).Value = 122200
$ws2.Cells.Item(37, 12).Value = 122200
$ws2.Cells.Item(37, 13).Value = 122200
$ws2.Cells.Item(37, 14).Value = 153000

$ws2.Cells.Item(38, 3).Value = 93200
$ws2.Cells.Item(38, 4).Value = 96200
$ws2.Cells.Item(38, 5).Value = 101200
$ws2.Cells.Item(38, 6).Value = 101200
$ws2.Cells.Item(38, 7).Value = 103200
$ws2.Cells.Item(38, 8).Value = 113200
$ws2.Cells.Item(38, 9).Value = 113200
$ws2.Cells.Item(38, 10).Value = 122200
$ws2.Cells.Item(38, 11).Value = 122200
$ws2.Cells.Item(38, 12).Value = 122200
$ws2.Cells.Item(38, 13).Value = 148000
$ws2.Cells.Item(38, 14).Value = 148000

# ---- Sheet2 (forecast): extend with new rows 39-50 (202104-202203) ----
$ws2.Range("B38:N38").Copy()
$ws2.Range("B39:N50").PasteSpecial(-4122)

# Row 39
$ws2.Cells.Item(39, 2).Value = 202104
$ws2.Cells.Item(39, 3).Value = 106000
$ws2.Cells.Item(39, 4).Value = 99600
$ws2.Cells.Item(39, 5).Value = 101200
$ws2.Cells.Item(39, 6).Value = 103200
$ws2.Cells.Item(39, 7).Value = 113200
$ws2.Cells.Item(39, 8).Value = 113200
$ws2.Cells.Item(39, 9).Value = 122200
$ws2.Cells.Item(39, 10).Value = 122200
$ws2.Cells.Item(39, 11).Value = 122200
$ws2.Cells.Item(39, 12).Value = 148000
$ws2.Cells.Item(39, 13).Value = 148000
$ws2.Cells.Item(39, 14).Value = 148000

# Row 40
$ws2.Cells.Item(40, 2).Value = 202105
$ws2.Cells.Item(40, 3).Value = 104200
$ws2.Cells.Item(40, 4).Value = 107300
$ws2.Cells.Item(40, 5).Value = 103620
$ws2.Cells.Item(40, 6).Value = 112640
$ws2.Cells.Item(40, 7).Value = 112640
$ws2.Cells.Item(40, 8).Value = 120060
$ws2.Cells.Item(40, 9).Value = 120340
$ws2.Cells.Item(40, 10).Value = 120620
$ws2.Cells.Item(40, 11).Value = 143000
$ws2.Cells.Item(40, 12).Value = 143000
$ws2.Cells.Item(40, 13).Value = 143000

# Row 41
$ws2.Cells.Item(41, 2).Value = 202106
$ws2.Cells.Item(41, 3).Value = 120320
$ws2.Cells.Item(41, 4).Value = 108000
$ws2.Cells.Item(41, 5).Value = 119940
$ws2.Cells.Item(41, 6).Value = 117540
$ws2.Cells.Item(41, 7).Value = 116760
$ws2.Cells.Item(41, 8).Value = 144100
$ws2.Cells.Item(41, 9).Value = 134220
$ws2.Cells.Item(41, 10).Value = 136600
$ws2.Cells.Item(41, 11).Value = 136600
$ws2.Cells.Item(41, 12).Value = 136600

# Row 42
$ws2.Cells.Item(42, 2).Value = 202107
$ws2.Cells.Item(42, 3).Value = 133200
$ws2.Cells.Item(42, 4).Value = 111940
$ws2.Cells.Item(42, 5).Value = 105940
$ws2.Cells.Item(42, 6).Value = 93960
$ws2.Cells.Item(42, 7).Value = 122900
$ws2.Cells.Item(42, 8).Value = 109040
$ws2.Cells.Item(42, 9).Value = 124200
$ws2.Cells.Item(42, 10).Value = 124200
$ws2.Cells.Item(42, 11).Value = 124200

# Row 43
$ws2.Cells.Item(43, 2).Value = 202108
$ws2.Cells.Item(43, 3).Value = 132140
$ws2.Cells.Item(43, 4).Value = 134880
$ws2.Cells.Item(43, 5).Value = 112800
$ws2.Cells.Item(43, 6).Value = 133200
$ws2.Cells.Item(43, 7).Value = 123600
$ws2.Cells.Item(43, 8).Value = 124200
$ws2.Cells.Item(43, 9).Value = 124200
$ws2.Cells.Item(43, 10).Value = 124200

# Row 44
$ws2.Cells.Item(44, 2).Value = 202109
$ws2.Cells.Item(44, 3).Value = 140400
$ws2.Cells.Item(44, 4).Value = 116600
$ws2.Cells.Item(44, 5).Value = 138079
$ws2.Cells.Item(44, 6).Value = 141435
$ws2.Cells.Item(44, 7).Value = 124200
$ws2.Cells.Item(44, 8).Value = 124200
$ws2.Cells.Item(44, 9).Value = 124200

# Row 45
$ws2.Cells.Item(45, 2).Value = 202110
$ws2.Cells.Item(45, 3).Value = 98800
$ws2.Cells.Item(45, 4).Value = 141400
$ws2.Cells.Item(45, 5).Value = 135200
$ws2.Cells.Item(45, 6).Value = 124200
$ws2.Cells.Item(45, 7).Value = 124200
$ws2.Cells.Item(45, 8).Value = 124200

# Row 46
$ws2.Cells.Item(46, 2).Value = 202111
$ws2.Cells.Item(46, 3).Value = 111200
$ws2.Cells.Item(46, 4).Value = 168800
$ws2.Cells.Item(46, 5).Value = 124200
$ws2.Cells.Item(46, 6).Value = 124200
$ws2.Cells.Item(46, 7).Value = 124200

# Row 47
$ws2.Cells.Item(47, 2).Value = 202112
$ws2.Cells.Item(47, 3).Value = 188600
$ws2.Cells.Item(47, 4).Value = 126000
$ws2.Cells.Item(47, 5).Value = 133249
$ws2.Cells.Item(47, 6).Value = 149016
$ws2.Cells.Item(47, 7).Value = 181375
$ws2.Cells.Item(47, 8).Value = 181375
$ws2.Cells.Item(47, 9).Value = 181375
$ws2.Cells.Item(47, 10).Value = 181375
$ws2.Cells.Item(47, 11).Value = 181375
$ws2.Cells.Item(47, 12).Value = 181375
$ws2.Cells.Item(47, 13).Value = 181375
$ws2.Cells.Item(47, 14).Value = 181375

# Row 48
$ws2.Cells.Item(48, 2).Value = 202201
$ws2.Cells.Item(48, 3).Value = 144400
$ws2.Cells.Item(48, 4).Value = 129400
$ws2.Cells.Item(48, 5).Value = 152800
$ws2.Cells.Item(48, 6).Value = 174888.45755693581
$ws2.Cells.Item(48, 7).Value = 174888.45755693581
$ws2.Cells.Item(48, 8).Value = 174888.45755693581
$ws2.Cells.Item(48, 9).Value = 174888.45755693581
$ws2.Cells.Item(48, 10).Value = 174888.45755693581
$ws2.Cells.Item(48, 11).Value = 174888.45755693581
$ws2.Cells.Item(48, 12).Value = 174888.45755693581
$ws2.Cells.Item(48, 13).Value = 174888.45755693581
$ws2.Cells.Item(48, 14).Value = 174888.45755693581

# Row 49
$ws2.Cells.Item(49, 2).Value = 202202
$ws2.Cells.Item(49, 3).Value = 105200
$ws2.Cells.Item(49, 4).Value = 143200
$ws2.Cells.Item(49, 5).Value = 174888.45755693581
$ws2.Cells.Item(49, 6).Value = 174888.45755693581
$ws2.Cells.Item(49, 7).Value = 174888.45755693581
$ws2.Cells.Item(49, 8).Value = 174888.45755693581
$ws2.Cells.Item(49, 9).Value = 174888.45755693581
$ws2.Cells.Item(49, 10).Value = 174888.45755693581
$ws2.Cells.Item(49, 11).Value = 174888.45755693581
$ws2.Cells.Item(49, 12).Value = 174888.45755693581
$ws2.Cells.Item(49, 13).Value = 174888.45755693581
$ws2.Cells.Item(49, 14).Value = 192500

# Row 50
$ws2.Cells.Item(50, 2).Value = 202203
$ws2.Cells.Item(50, 3).Value = 102000
$ws2.Cells.Item(50, 4).Value = 151500
$ws2.Cells.Item(50, 5).Value = 148000
$ws2.Cells.Item(50, 6).Value = 148000
$ws2.Cells.Item(50, 7).Value = 172000
$ws2.Cells.Item(50, 8).Value = 156000
$ws2.Cells.Item(50, 9).Value = 156000
$ws2.Cells.Item(50, 10).Value = 164000
$ws2.Cells.Item(50, 11).Value = 164000
$ws2.Cells.Item(50, 12).Value = 164000
$ws2.Cells.Item(50, 13).Value = 177013.38199513382
$ws2.Cells.Item(50, 14).Value = 177013.38199513382

# ---- Sheet view state ----
# Reflect the final view: sheet2 scrolled/selected near its new last row (B50),
# while sheet1 (the first/actuals tab) remains the active sheet.
$ws2.Activate()
$ws2.Range("B50").Select()
$ws1.Activate()

Write-Output "edit complete"
